# 8.7.1 workbook update
# Refreshes the "by sex" / "urban-rural" / a few English & Kyrgyz labels on
# the single data sheet so the RU/EN/KY translation triples line up again
# (e.g. "Мужской"/"Male " -> "Мужчины"/"Men", "шаар"/"urban" -> "Шаар"/"Urban",
# capitalised English labels in the education/wealth-quintile block, and a
# previously-blank Kyrgyz header for "Functional difficulties in a child").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: "by sex" section header
$ws.Range("A6").Value = "Жынысы боюнча"
$ws.Range("B6").Value = "По полу"
$ws.Range("C6").Value = "By sex"

# Row 7: male
$ws.Range("A7").Value = "Эркектер"
$ws.Range("B7").Value = "Мужчины"
$ws.Range("C7").Value = "Men"

# Row 8: female
$ws.Range("A8").Value = "Аялдар"
$ws.Range("B8").Value = "Женщины"
$ws.Range("C8").Value = "Woman"

# Row 10: urban
$ws.Range("A10").Value = "Шаар"
$ws.Range("C10").Value = "Urban"

# Row 11: rural
$ws.Range("A11").Value = "Айыл"
$ws.Range("C11").Value = "Rural"

# Row 28: does not attend (school attendance, English label)
$ws.Range("C28").Value = "Does not attend"

# Row 29: education of mother (English label)
$ws.Range("C29").Value = "Educationof mother"

# Row 30: preschool or not / primary (English label)
$ws.Range("C30").Value = "Preschool or not /primary"

# Row 31: basic general (English label)
$ws.Range("C31").Value = "Basic general"

# Row 32: average total (English label)
$ws.Range("C32").Value = "Average total"

# Row 33: vocational primary/secondary (English label)
$ws.Range("C33").Value = "Vocational primary /secondary"

# Row 34: higher (English label)
$ws.Range("C34").Value = "Higher"

# Row 35: functional difficulties in a child - add the missing Kyrgyz label
$ws.Range("A35").Value = "Баланын функционалдык кыйнчылыктары"

# Row 38: wealth quintile (English label)
$ws.Range("C38").Value = "Wealth quintile"
